$wb = $excel.ActiveWorkbook

$wsStart = $wb.Worksheets.Item("Ausgangstabelle")
$wsMon   = $wb.Worksheets.Item("InMonateZerlegen")

# --- Value edits: demo date changed from 2023-03-31 (45016) to 2023-03-20 (45005) ---
$wsStart.Range("C5").Value = 45005
$wsMon.Range("C29:D40").Value = 45005

# --- Selection / active sheet changes ---
# EndeBestimmnen keeps its own selection untouched.

# Select C3 on Ausgangstabelle (but don't leave it the active sheet yet)
$wsStart.Activate()
$wsStart.Range("C3").Select() | Out-Null

# Select D27 on InMonateZerlegen
$wsMon.Activate()
$wsMon.Range("D27").Select() | Out-Null

# Final active sheet is Ausgangstabelle (tab 0), matching removal of activeTab="2"
$wsStart.Activate()
